$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'" + '71.004.86'
$ws.Range('E2').Value = '  -0.25%  '

$ws.Range('D3').Value = "'" + '3.860.45'
$ws.Range('E3').Value = '  +1.35%  '

$ws.Range('D4').Value = "'" + '0.999'
$ws.Range('E4').Value = '  -0.08%  '

$ws.Range('D5').Value = "'" + '701.15'
$ws.Range('E5').Value = '  +0.37%  '

$ws.Range('D6').Value = "'" + '173.33'
$ws.Range('E6').Value = '  +0.20%  '

$ws.Range('D7').Value = "'" + '3.859.58'
$ws.Range('E7').Value = '  +1.40%  '

$ws.Range('E8').Value = '  +0.03%  '

$ws.Range('E9').Value = '  -0.14%  '

$ws.Range('E10').Value = '  -0.45%  '

$ws.Range('D11').Value = "'" + '7.22'
$ws.Range('E11').Value = '  -3.79%  '

$ws.Range('E12').Value = '  -0.54%  '

$ws.Range('E13').Value = '  -0.19%  '

$ws.Range('E14').Value = '  -0.05%  '

$ws.Range('D15').Value = "'" + '4.510.91'
$ws.Range('E15').Value = '  +1.39%  '

$ws.Range('D16').Value = "'" + '3.971.57'
$ws.Range('E16').Value = '  +4.34%  '

$ws.Range('D17').Value = "'" + '71.055.95'
$ws.Range('E17').Value = '  -0.12%  '

$ws.Range('D18').Value = "'" + '7.21'
$ws.Range('E18').Value = '  -0.24%  '

$ws.Range('D19').Value = "'" + '17.43'
$ws.Range('E19').Value = '  -2.55%  '

$ws.Range('E20').Value = '  -0.47%  '

$ws.Range('D21').Value = "'" + '501.56'
$ws.Range('E21').Value = '  +4.21%  '

$ws.Range('D22').Value = "'" + '10.74'
$ws.Range('E22').Value = '  -4.14%  '

$ws.Range('D23').Value = "'" + '0.721'
$ws.Range('E23').Value = '  +0.77%  '

$ws.Range('D24').Value = "'" + '0.0000149'
$ws.Range('E24').Value = '  +2.95%  '

$ws.Range('D25').Value = "'" + '85.04'
$ws.Range('E25').Value = '  +1.23%  '

$ws.Range('D26').Value = "'" + '10.66'
$ws.Range('E26').Value = '  +1.60%  '

$ws.Range('D27').Value = "'" + '12.22'
$ws.Range('E27').Value = '  -1.41%  '

$ws.Range('E28').Value = '  -2.15%  '

$ws.Range('D29').Value = "'" + '3.16'
$ws.Range('E29').Value = '  +2.62%  '

$ws.Range('D31').Value = "'" + '7.55'
$ws.Range('E31').Value = '  -0.11%  '

$ws.Range('E32').Value = '  -1.88%  '

$ws.Range('B33').Value = 'EthereumClassic'
$ws.Range('C33').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D33').Value = "'" + '29.59'
$ws.Range('E33').Value = '  +0.04%  '

$ws.Range('B34').Value = 'Kaspa'
$ws.Range('C34').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D34').Value = "'" + '0.182'
$ws.Range('E34').Value = '  +3.60%  '

$ws.Range('D35').Value = "'" + '9.22'
$ws.Range('E35').Value = '  -0.14%  '

$ws.Range('D36').Value = "'" + '3.816.99'
$ws.Range('E36').Value = '  +1.53%  '

$ws.Range('D37').Value = "'" + '0.999'
$ws.Range('E37').Value = '  -0.05%  '

$ws.Range('E38').Value = '  +1.23%  '

$ws.Range('D39').Value = "'" + '2.39'
$ws.Range('E39').Value = '  +7.74%  '

$ws.Range('E40').Value = '  +8.67%  '

$ws.Range('E41').Value = '  -1.80%  '

$ws.Range('D42').Value = "'" + '6.04'
$ws.Range('E42').Value = '  +0.88%  '

$ws.Range('E44').Value = '  +0.15%  '

$ws.Range('E45').Value = '  -3.15%  '

$ws.Range('E46').Value = '  +2.14%  '

$ws.Range('D47').Value = "'" + '49.23'
$ws.Range('E47').Value = '  -0.09%  '

$ws.Range('D48').Value = "'" + '417.68'
$ws.Range('E48').Value = '  +2.96%  '

$ws.Range('E49').Value = '  +0.43%  '

$ws.Range('E50').Value = '  -2.37%  '

$ws.Range('D51').Value = "'" + '43.42'
$ws.Range('E51').Value = '  -4.65%  '
